# [MOSIP-14369] Fix: boolean values
#
# The E column (is_active) previously stored its values as the formula
# =TRUE() (a numeric boolean result). This rewrites E2:E11 so they hold
# the literal text "TRUE" instead (a plain string, not a boolean),
# matching the corrected workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    # Writing the bare string "TRUE" via .Value gets auto-coerced to a
    # boolean by Excel's normal input parsing (same as typing TRUE into a
    # cell). Routing it through a text formula and collapsing the formula
    # to its computed value via copy / paste-values keeps the cell typed
    # as text instead.
    $cell.Formula = "=""TRUE"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

$ws.Range("E2:E11").Select()
